$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 847, shifting existing rows 847:888 down to 848:889
$ws.Range("A847").EntireRow.Insert()

# Populate the newly inserted row 847 with the new data point.
# Column A holds a date-like string ("2026/02/22") that must stay a literal
# text value (matching the rest of the sheet, which stores dates as plain
# text rather than real date serials). Force text entry with NumberFormat="@"
# so Excel doesn't auto-convert it to a date serial, then restore the
# cell's style to Normal so no stray formatting is left behind.
$ws.Range("A847").NumberFormat = "@"
$ws.Range("A847").Value = "2026/02/22"
$ws.Range("A847").NumberFormat = "General"
$ws.Range("A847").Style = "Normal"

$ws.Range("B847").Value = "日"
$ws.Range("C847").Value = 5
$ws.Range("D847").Value = 201
